$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 266; $r++) {
    $ws.Cells.Item($r, 2).Value = 2023
}

$ws.Cells.Item(1211, 2).Value = 45488
$ws.Cells.Item(1277, 2).Value = 45488
$ws.Cells.Item(1298, 2).Value = 45488
